$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H20").Value = 1
$ws.Range("H21").Value = $null
$ws.Range("H24").Value = 1
$ws.Range("H27").Value = 1

$ws.Range("H28").Select() | Out-Null
